$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# December 2020 (Dez / column Q) production figures, rows 2-93
$decValues = @(
    708,1479,139,3178,583,0,3,4,20312,14849,
    14027,38645,9646,820,221,20257,38,273,698,5179,
    1649,688,13828,12331,32,33458,1636,395,0,19426,
    2404,3569,695,288,0,313,658,957,281,184,
    1405,3306,4218,204,2,283,0,68,54,220,
    9,0,2,3,30,1,104,319,533,13,
    2,0,2,56,484,183,0,0,0,5,
    234,1323,891,2,230,2,286,0,131,243,
    0,1278,1319,1525,1,169,137,372,359,30,
    52,26
)

$startRow = 2
for ($i = 0; $i -lt $decValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 17).Value = $decValues[$i]
}
